# Adds five new "semaforo di gruppo" indicator rows to the Library_Formula
# sheet (INDICATOR_QUOTA_EXP_GRP1/2/3, CALCOLO_MODULO, SEMAFORO_DI_GRUPPO),
# following the existing CREATE/MODIFY | LIB_EWS_IT | <name> | | String | String
# pattern used by all the other indicator rows, and updates the sheet
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp
$newCount = 5

# Append new rows by duplicating the last existing data row (keeps the
# same styling / column layout as the rest of the indicator list), then
# fill in column C for each of the new rows.
for ($i = 0; $i -lt $newCount; $i++) {
    $ws.Rows.Item($lastRow).Copy()
    $ws.Rows.Item($lastRow + 1).Insert(-4121)   # xlShiftDown
    $lastRow = $lastRow + 1
}

$firstNewRow = $lastRow - $newCount + 1

$ws.Cells.Item($firstNewRow, 3).Value     = "INDICATOR_QUOTA_EXP_GRP1"
$ws.Cells.Item($firstNewRow + 1, 3).Value = "INDICATOR_QUOTA_EXP_GRP2"
$ws.Cells.Item($firstNewRow + 2, 3).Value = "INDICATOR_QUOTA_EXP_GRP3"
$ws.Cells.Item($firstNewRow + 4, 3).Value = "SEMAFORO_DI_GRUPPO"
$ws.Cells.Item($firstNewRow + 3, 3).Value = "CALCOLO_MODULO"

# Update the current selection on the sheet.
$ws.Range("E2").Select()
